# Fix Training Data Issue (#48)
# The "Date" column (BF) was stamped with the source filename-derived
# string "6-10-2007-08" for every data row. NBA.com stats for a given
# "game date" are actually posted the day *after* that date, so the
# correct calendar date for this file is 2008-06-10. Correct every
# row's BF cell (BF2:BF31) to the proper ISO date text "2008-06-10".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 31
$col = 58   # column BF
$correctedDate = "2008-06-10"

# Format the column as text first so Excel stores the corrected value
# as a literal string instead of re-interpreting "2008-06-10" as a
# date serial (it would otherwise silently roll it back to a date
# value because the text looks like an ISO date).
$ws.Range("BF$firstRow`:BF$lastRow").NumberFormat = "@"

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, $col).Value = $correctedDate
}
